$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'isophonics_232'
$ws.Range("B2").Value = 'isophonics_261'
$ws.Range("C2").Value = 0.1770833333333333
$ws.Range("D2").Value = '[[''C'', ''F'', ''C'']]'
$ws.Range("E2").Value = '[[''G/3'', ''C'', ''G'']]'
$ws.Range("F2").Value = '[(16.17907, 21.844739)]'
$ws.Range("G2").Value = '[(13.079206, 17.500916)]'
$ws.Range("H2").Value = ''
$ws.Range("I2").Value = 'spotify:track:2B4Y9u4ERAFiMo13XPJyGP'

# Row 3
$ws.Range("A3").Value = 'schubert-winterreise_123'
$ws.Range("B3").Value = 'schubert-winterreise_160'
$ws.Range("C3").Value = 0.2701149425287356
$ws.Range("D3").Value = '[[''F:min'', ''C:maj/F'', ''F:min'', ''C:maj/F'', ''F:min'', ''C:maj/F'']]'
$ws.Range("E3").Value = '[[''D:min'', ''A:maj'', ''D:min'', ''A:maj'', ''D:min'', ''A:maj'']]'
$ws.Range("F3").Value = '[(21.9, 51.5)]'
$ws.Range("G3").Value = '[(10.82, 26.2)]'
$ws.Range("H3").Value = ''
$ws.Range("I3").Value = ''

# Row 4
$ws.Range("A4").Value = 'schubert-winterreise_101'
$ws.Range("B4").Value = 'schubert-winterreise_14'
$ws.Range("C4").Value = 0.3333333333333333
$ws.Range("D4").Value = '[[''F#:7/B'', ''F#:(3,5,b7,b9)/B'', ''B:min'']]'
$ws.Range("E4").Value = '[[''A:7'', ''A:(3,5,b7,b9)/G'', ''D:min/F'']]'
$ws.Range("F4").Value = '[(4.7, 21.62)]'
$ws.Range("G4").Value = '[(165.14, 168.7)]'
$ws.Range("H4").Value = ''
$ws.Range("I4").Value = ''

# Row 5
$ws.Range("A5").Value = 'isophonics_160'
$ws.Range("B5").Value = 'schubert-winterreise_173'
$ws.Range("C5").Value = 0.1169415292353823
$ws.Range("D5").Value = '[[''Eb:7'', ''Ab:maj'', ''Eb/3''], [''Eb:maj'', ''Bb:7/3'', ''Eb:maj''], [''Eb:7'', ''Ab:maj'', ''Ab/7'']]'
$ws.Range("E5").Value = '[[''D:7'', ''G:maj'', ''D:maj''], [''D:maj/A'', ''A:7'', ''D:maj''], [''D:7'', ''G:maj'', ''G:maj/B'']]'
$ws.Range("F5").Value = '[(26.006, 29.321), (42.588, 46.723), (54.072, 59.053)]'
$ws.Range("G5").Value = '[(42.62, 45.06), (8.7, 10.94), (18.3, 20.68)]'
$ws.Range("H5").Value = ''
$ws.Range("I5").Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

# Row 6
$ws.Range("A6").Value = 'jaah_43'
$ws.Range("B6").Value = 'isophonics_235'
$ws.Range("C6").Value = 0.119047619047619
$ws.Range("D6").Value = '[[''Eb'', ''Eb:7'', ''Ab'']]'
$ws.Range("E6").Value = '[[''Bb'', ''Bb:7'', ''Eb'']]'
$ws.Range("F6").Value = '[(45.01, 49.44)]'
$ws.Range("G6").Value = '[(48.170328, 55.589104)]'
$ws.Range("H6").Value = ''
$ws.Range("I6").Value = ''

# Row 7
$ws.Range("A7").Value = 'isophonics_208'
$ws.Range("B7").Value = 'isophonics_2'
$ws.Range("C7").Value = 0.2371794871794872
$ws.Range("D7").Value = '[[''E'', ''B'', ''E'', ''B'']]'
$ws.Range("E7").Value = '[[''Eb:maj'', ''Bb:maj'', ''Eb:maj'', ''Bb:maj'']]'
$ws.Range("F7").Value = '[(148.857891, 156.114126)]'
$ws.Range("G7").Value = '[(217.911, 223.713)]'
$ws.Range("H7").Value = ''
$ws.Range("I7").Value = ''

# Row 8
$ws.Range("A8").Value = 'isophonics_109'
$ws.Range("B8").Value = 'isophonics_152'
$ws.Range("C8").Value = 0.5576923076923077
$ws.Range("D8").Value = '[[''A'', ''D'', ''A'']]'
$ws.Range("E8").Value = '[[''A/3'', ''D'', ''A'']]'
$ws.Range("F8").Value = '[(94.925, 100.821)]'
$ws.Range("G8").Value = '[(4.25, 9.179)]'
$ws.Range("H8").Value = ''
$ws.Range("I8").Value = ''

# Row 9
$ws.Range("A9").Value = 'isophonics_81'
$ws.Range("B9").Value = 'jaah_39'
$ws.Range("C9").Value = 0.1940639269406393
$ws.Range("D9").Value = '[[''E/5'', ''B'', ''B'', ''E/5'']]'
$ws.Range("E9").Value = '[[''Db'', ''Ab'', ''Ab'', ''Db'']]'
$ws.Range("F9").Value = '[(61.620521, 70.444104)]'
$ws.Range("G9").Value = '[(126.19, 129.59)]'
$ws.Range("H9").Value = ''
$ws.Range("I9").Value = ''

# Row 10
$ws.Range("A10").Value = 'isophonics_152'
$ws.Range("B10").Value = 'isophonics_273'
$ws.Range("C10").Value = 0.1031468531468532
$ws.Range("D10").Value = '[[''A/3'', ''E/5'', ''A/3'']]'
$ws.Range("E10").Value = '[[''C'', ''G'', ''C'']]'
$ws.Range("F10").Value = '[(14.168, 18.573)]'
$ws.Range("G10").Value = '[(20.801, 27.165)]'
$ws.Range("H10").Value = ''
$ws.Range("I10").Value = ''

# Row 11
$ws.Range("A11").Value = 'schubert-winterreise_68'
$ws.Range("B11").Value = 'schubert-winterreise_172'
$ws.Range("C11").Value = 0.7083333333333333
$ws.Range("D11").Value = '[[''D:7'', ''G:maj'', ''D:7'', ''G:maj'', ''D:7'', ''G:maj'']]'
$ws.Range("E11").Value = '[[''B:7'', ''E:maj'', ''B:7'', ''E:maj'', ''B:7'', ''E:maj'']]'
$ws.Range("F11").Value = '[(33.0, 48.22)]'
$ws.Range("G11").Value = '[(14.44, 41.8)]'
$ws.Range("H11").Value = ''
$ws.Range("I11").Value = ''

# Row 12
$ws.Range("A12").Value = 'schubert-winterreise_13'
$ws.Range("B12").Value = 'schubert-winterreise_33'
$ws.Range("C12").Value = 0.2767857142857143
$ws.Range("D12").Value = '[[''B:maj'', ''F#:7'', ''B:maj''], [''F#:7'', ''B:maj'', ''B:min'']]'
$ws.Range("E12").Value = '[[''G:maj/D'', ''D:7'', ''G:maj''], [''D:7'', ''G:maj'', ''G:min'']]'
$ws.Range("F12").Value = '[(0.32, 9.54), (30.52, 40.64)]'
$ws.Range("G12").Value = '[(65.44, 66.82), (66.04, 69.18)]'
$ws.Range("H12").Value = ''
$ws.Range("I12").Value = ''

# Row 13
$ws.Range("A13").Value = 'schubert-winterreise_186'
$ws.Range("B13").Value = 'schubert-winterreise_123'
$ws.Range("C13").Value = 0.2528735632183908
$ws.Range("D13").Value = '[[''F:min'', ''C:7'', ''F:min'', ''C:maj'', ''F:min'']]'
$ws.Range("E13").Value = '[[''F:min'', ''C:7/F'', ''F:min'', ''C:maj/F'', ''F:min'']]'
$ws.Range("F13").Value = '[(11.74, 19.74)]'
$ws.Range("G13").Value = '[(18.82, 46.24)]'
$ws.Range("H13").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I13").Value = ''

# Row 14
$ws.Range("A14").Value = 'schubert-winterreise_44'
$ws.Range("B14").Value = 'schubert-winterreise_91'
$ws.Range("C14").Value = 0.09642857142857142
$ws.Range("D14").Value = '[[''A#/F'', ''F:7'', ''A#'']]'
$ws.Range("E14").Value = '[[''D#:maj/A#'', ''A#:7'', ''D#:maj'']]'
$ws.Range("F14").Value = '[(271.22, 275.78)]'
$ws.Range("G14").Value = '[(79.98, 87.7)]'
$ws.Range("H14").Value = ''
$ws.Range("I14").Value = ''

# Row 15
$ws.Range("A15").Value = 'isophonics_128'
$ws.Range("B15").Value = 'schubert-winterreise_199'
$ws.Range("C15").Value = 0.2932330827067669
$ws.Range("D15").Value = '[[''C'', ''F'', ''C'']]'
$ws.Range("E15").Value = '[[''G:maj/B'', ''C:maj'', ''G:maj/D'']]'
$ws.Range("F15").Value = '[(10.634761, 15.986961)]'
$ws.Range("G15").Value = '[(62.2, 64.4)]'
$ws.Range("H15").Value = ''
$ws.Range("I15").Value = ''

# Row 16
$ws.Range("A16").Value = 'schubert-winterreise_108'
$ws.Range("B16").Value = 'schubert-winterreise_63'
$ws.Range("C16").Value = 0.1666666666666667
$ws.Range("D16").Value = '[[''A:min/E'', ''E:7'', ''A:min'', ''A:min/E''], [''C:maj'', ''A:min/E'', ''E:7'', ''A:min''], [''A:min'', ''A:min/E'', ''E:7'', ''A:min'']]'
$ws.Range("E16").Value = '[[''C:min/G'', ''G:7'', ''C:min'', ''C:min''], [''D#/G'', ''C:min/G'', ''G:7'', ''C:min''], [''C:min'', ''C:min/G'', ''G:7'', ''C:min'']]'
$ws.Range("F16").Value = '[(32.7, 39.58), (30.76, 37.72), (36.02, 45.64)]'
$ws.Range("G16").Value = '[(25.8, 31.82), (67.16, 72.4), (24.68, 29.46)]'
$ws.Range("H16").Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'
$ws.Range("I16").Value = ''

# Row 17
$ws.Range("A17").Value = 'isophonics_204'
$ws.Range("B17").Value = 'isophonics_203'
$ws.Range("C17").Value = 0.1045673076923077
$ws.Range("D17").Value = '[[''D'', ''A'', ''A'']]'
$ws.Range("E17").Value = '[[''G'', ''D'', ''D/b7'']]'
$ws.Range("F17").Value = '[(0.459543, 6.438299)]'
$ws.Range("G17").Value = '[(55.397573, 60.285374)]'
$ws.Range("H17").Value = ''
$ws.Range("I17").Value = ''
